# Add the "Assignment_17" row to the assignments tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row values
$ws.Range("A18").Value = "Assignment_17"
$ws.Range("B18").Value = "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_17/MovieAPI"
$ws.Range("C18").Value = 45173

# Add the hyperlink for the new assignment's GitHub link cell
$ws.Hyperlinks.Add($ws.Range("B18"), "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_17/MovieAPI")

# Mirror the formatting of the previous row (A17:C17 -> A18:C18) so the new
# row picks up the same styles (border, font, hyperlink look, date number
# format, row height, etc.) - applied AFTER the hyperlink so the paste wins
# over Excel's automatic "Hyperlink" style stamp.
$ws.Range("A17:C17").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(18).RowHeight = 34.5

# Update the selection to match where the author ended up after entering the row
$ws.Range("B18:C18").Select() | Out-Null
